$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.272.08'
$ws.Range("E2").Value = '  +1.17%  '
$ws.Range("D3").Value = '1.566.85'
$ws.Range("E3").Value = '  +0.54%  '
$ws.Range("E4").Value = '  -0.13%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '210.86'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.38%  '
$ws.Range("E6").Value = '  +0.11%  '
$ws.Range("E7").Value = '  -0.29%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '22.25'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.66%  '
$ws.Range("E9").Value = '  +0.28%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0597'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.29%  '
$ws.Range("E11").Value = '  +1.72%  '
$ws.Range("D12").Value = '1.789.57'
$ws.Range("E12").Value = '  +0.50%  '
$ws.Range("D13").Value = '1.567.32'
$ws.Range("E13").Value = '  +1.89%  '
$ws.Range("E14").Value = '  +0.23%  '
$ws.Range("E15").Value = '  +0.01%  '
$ws.Range("D16").Value = '27.222.22'
$ws.Range("E16").Value = '  +0.98%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '62.00'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.14%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '7.48'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.57%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '217.55'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.73%  '
$ws.Range("D20").Value = '0.0₃0702'
$ws.Range("E21").Value = '  -0.21%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.15'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.84%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.34'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.27%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.96'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.50%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '151.54'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.95%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.63'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.50%  '
$ws.Range("E27").Value = '  +1.97%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.05'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.26%  '
$ws.Range("E29").Value = '  -0.27%  '
$ws.Range("E30").Value = '  +1.57%  '
$ws.Range("E31").Value = '  -0.31%  '
$ws.Range("E32").Value = '  +0.36%  '
$ws.Range("D33").Value = '1.459.10'
$ws.Range("E33").Value = '  +2.53%  '
$ws.Range("E34").Value = '  +0.26%  '
$ws.Range("E35").Value = '  +4.45%  '
$ws.Range("E36").Value = '  +1.63%  '
$ws.Range("E37").Value = '  +0.39%  '
$ws.Range("E38").Value = '  +0.05%  '
$ws.Range("E39").Value = '  +1.56%  '
$ws.Range("E40").Value = '  +0.66%  '
$ws.Range("E41").Value = '  +0.98%  '
$ws.Range("E42").Value = '  -0.23%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.35'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.44%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.987'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.78%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '64.54'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.01%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.75'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.17%  '
$ws.Range("D47").Value = '1.701.40'
$ws.Range("E47").Value = '  +0.46%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '85.98'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.35%  '
$ws.Range("B49").Value = 'Cronos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0525'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.17%  '
$ws.Range("B50").Value = 'BabyDogeCoin'
$ws.Range("C50").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D50").Value = '0.0₆0100'
$ws.Range("E50").Value = '  -1.78%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0949'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.08%  '
